# Reorders data rows 2-14 on the active sheet according to a fixed row
# permutation (observation records were re-sorted upstream; this replays
# the same reordering against the local copy of the export).
#
# Row r=11 is unchanged; the rest map as:
#   new row 2  <- old row 5
#   new row 3  <- old row 6
#   new row 4  <- old row 14
#   new row 5  <- old row 7
#   new row 6  <- old row 8
#   new row 7  <- old row 3
#   new row 8  <- old row 9
#   new row 9  <- old row 10
#   new row 10 <- old row 2
#   new row 11 <- old row 11
#   new row 12 <- old row 13
#   new row 13 <- old row 12
#   new row 14 <- old row 4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 14
$numRows = $lastRow - $firstRow + 1
$firstCol = 1
$lastCol = 51   # column AY
$numCols = $lastCol - $firstCol + 1

# new row (key) -> old row (value), both in worksheet row numbers
$rowMap = @{
    2  = 5
    3  = 6
    4  = 14
    5  = 7
    6  = 8
    7  = 3
    8  = 9
    9  = 10
    10 = 2
    11 = 11
    12 = 13
    13 = 12
    14 = 4
}

# Columns that hold genuine numbers / booleans in this export; every other
# column is free text (inlineStr) and must round-trip as text even when its
# content looks like a number or a date, so we force those columns to a
# Text number format before writing the permuted values back.
$numericCols = @(1, 2, 5, 17, 18, 19)        # A, B, E, Q, R, S
$booleanCols = @(30, 31, 33)                  # AD, AE, AG

$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$data = $srcRange.Value2

# Force every non-numeric, non-boolean column to Text format so the
# write-back below can't get reinterpreted as a number or a date.
for ($c = $firstCol; $c -le $lastCol; $c++) {
    if ($numericCols -notcontains $c -and $booleanCols -notcontains $c) {
        $colRange = $ws.Range($ws.Cells.Item($firstRow, $c), $ws.Cells.Item($lastRow, $c))
        $colRange.NumberFormat = "@"
    }
}

$newData = New-Object 'object[,]' $numRows, $numCols
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    $srcIdx = $oldRow - $firstRow + 1
    $dstIdx = $newRow - $firstRow
    for ($c = 1; $c -le $numCols; $c++) {
        $newData[$dstIdx, $c - 1] = $data[$srcIdx, $c]
    }
}

$srcRange.Value2 = $newData
